# Hortaliza, Vega Modelo de Temuco - Cebollín : weekly data refresh.
# A new weekly record is inserted at row 159 (the existing rows 159-207
# shift down to 160-208); the sheet's used range grows from A1:R207 to
# A1:R208.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 159 - this shifts rows
# 159..207 down to 160..208 and preserves their formatting/content.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new reading.
$ws.Range('A159').Value = 10
$ws.Range('B159').Value = 'Vega Modelo de Temuco'
$ws.Range('C159').Value = 'La Araucanía'
$ws.Range('D159').Value = 44463
$ws.Range('E159').Value = 9
$ws.Range('F159').Value = 100112037
$ws.Range('G159').Value = 'Cebollín'
$ws.Range('H159').Value = 'Sin especificar'
$ws.Range('I159').Value = 'Primera'
$ws.Range('J159').Value = 20
$ws.Range('K159').Value = 8000
$ws.Range('L159').Value = 8000
$ws.Range('M159').Value = 8000
$ws.Range('N159').Value = '$/docena de paquetes'
$ws.Range('O159').Value = 'Provincia de Cautín'
$ws.Range('P159').Value = 667
$ws.Range('Q159').Value = 12
$ws.Range('R159').Value = 'Hortaliza'
